# Regenerate the per-task-order worksheets with a fresh randomized
# experiment order (new script run), and re-order the "NB" task-order
# sheet so it now follows "TOL" in the tab order:
#   before: GNG, NB, RS, TOL, vSAT
#   after : GNG, RS, TOL, NB, vSAT

$wb = $excel.ActiveWorkbook

# --- Reorder sheets -------------------------------------------------
# Move the NB sheet (currently 2nd) so it sits right after TOL (4th).
$wsNB  = $wb.Worksheets.Item("NB_TO-1651255540253471")
$wsTOL = $wb.Worksheets.Item("TOL_TO-16512555403014758")
$wsNB.Move($null, $wsTOL)

# --- Rename sheets to the new generated task-order ids ---------------
$wsGNG  = $wb.Worksheets.Item("GNG_TO-16512555378532143")
$wsRS   = $wb.Worksheets.Item("RS_TO-16512555402544723")
$wsTOL  = $wb.Worksheets.Item("TOL_TO-16512555403014758")
$wsNB   = $wb.Worksheets.Item("NB_TO-1651255540253471")
$wsvSAT = $wb.Worksheets.Item("vSAT_TO-1651255540380471")

$wsGNG.Name  = "GNG_TO-1651588970469845"
$wsRS.Name   = "RS_TO-16515889704728415"
$wsTOL.Name  = "TOL_TO-16515889705326967"
$wsNB.Name   = "NB_TO-16515889735058324"
$wsvSAT.Name = "vSAT_TO-1651588973585828"

# --- Update GNG stim order -------------------------------------------
$wsGNG.Range("B2").Value = "go_stims-16515889704266326.csv"
$wsGNG.Range("B3").Value = "GNG_stims-1651588970443203.csv"
$wsGNG.Range("B4").Value = "go_stims-16515889704441674.csv"
$wsGNG.Range("B5").Value = "GNG_stims-16515889704688404.csv"

# --- RS sheet (eyes closed / eyes open) stays the same ----------------
$wsRS.Range("B2").Value = "eyes closed"
$wsRS.Range("B3").Value = "eyes open"

# --- Update TOL stim order (MM / ZM stims) ----------------------------
$wsTOL.Range("B2").Value = "MM_stims-16515889705006914.csv"
$wsTOL.Range("B3").Value = "ZM_stims-1651588970476392.csv"
$wsTOL.Range("B4").Value = "MM_stims-16515889705156953.csv"
$wsTOL.Range("B5").Value = "ZM_stims-16515889705026975.csv"
$wsTOL.Range("B6").Value = "MM_stims-16515889705316956.csv"
$wsTOL.Range("B7").Value = "ZM_stims-16515889705166957.csv"

# --- Update NB stim order (OB / ZB-match / TB stims) -------------------
$wsNB.Range("B2").Value  = "OB-16515889727378304.csv"
$wsNB.Range("B3").Value  = "ZB-match_4-1651588970669989.csv"
$wsNB.Range("B4").Value  = "OB-16515889712966251.csv"
$wsNB.Range("B5").Value  = "TB-16515889734258292.csv"
$wsNB.Range("B6").Value  = "ZB-match_7-165158897074455.csv"
$wsNB.Range("B7").Value  = "OB-1651588970959285.csv"
$wsNB.Range("B8").Value  = "TB-16515889731488318.csv"
$wsNB.Range("B9").Value  = "TB-16515889734938266.csv"
$wsNB.Range("B10").Value = "ZB-match_2-16515889705668998.csv"

# --- Update vSAT stim order --------------------------------------------
$wsvSAT.Range("B2").Value = "vSAT_stims-16515889735698276.csv"
$wsvSAT.Range("B3").Value = "SAT_stims-165158897351383.csv"
$wsvSAT.Range("B4").Value = "vSAT_stims-16515889735538316.csv"
$wsvSAT.Range("B5").Value = "SAT_stims-16515889735378273.csv"
